# Bugfix: remove special characters so generated .NET property names are valid.
# Adds a set of "Prop*" sample columns (G:L) to Plan1 exercising different
# special characters (space, backslash, slash, dash, underscore) so the
# ExcelProvider type-provider tests can confirm the sanitizer strips them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("G1").Value = "Prop sample"
$ws.Range("L1").Value = "Prop  \  sample  /  a   b"
$ws.Range("G2").Value = "Teste"
$ws.Range("H1").Value = "Prop-sample 2"
$ws.Range("I1").Value = "Prop_sample 3"
$ws.Range("J1").Value = "Prop\sample 4"
$ws.Range("K1").Value = "Prop/sample 5"

# --- Data rows (rows 2-7), columns G:L all filled with "Teste" ----------
$ws.Range("H2").Value = "Teste"
$ws.Range("I2").Value = "Teste"
$ws.Range("J2").Value = "Teste"
$ws.Range("K2").Value = "Teste"
$ws.Range("L2").Value = "Teste"

for ($row = 3; $row -le 7; $row++) {
    $ws.Cells.Item($row, 7).Value = "Teste"
    $ws.Cells.Item($row, 8).Value = "Teste"
    $ws.Cells.Item($row, 9).Value = "Teste"
    $ws.Cells.Item($row, 10).Value = "Teste"
    $ws.Cells.Item($row, 11).Value = "Teste"
    $ws.Cells.Item($row, 12).Value = "Teste"
}

# --- Column widths: mirror the auto-fit result Excel computed for these
# new columns (bestFit). ColumnWidth is expressed in character units; the
# saved <col> width = ColumnWidth + 0.8333333333333334.
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666   # -> width 12
$ws.Columns.Item(8).ColumnWidth = 13.0                 # -> width 13.833333333333334 (~13.85546875)
$ws.Columns.Item(9).ColumnWidth = 13.333333333333332   # -> width 14.166666666666666 (~14.140625)
$ws.Columns.Item(10).ColumnWidth = 13.166666666666666  # -> width 14
$ws.Columns.Item(11).ColumnWidth = 13.166666666666666  # -> width 14
$ws.Columns.Item(12).ColumnWidth = 19.666666666666668  # -> width 20.5 (~20.42578125)

# --- Selection left on L4, matching the saved view ----------------------
$ws.Range("L4").Select() | Out-Null
